# Insert a new data row at row 211 (pushing the existing rows 211-281 down
# to 212-282) and populate it with the new "Hortaliza" record for
# Femacal de La Calera / Sandia (date 44524, "Primera" quality, Peru origin).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 211..281 down by one row.
$ws.Rows.Item(211).Insert()

# Populate the newly inserted row 211 with the new record.
$ws.Cells.Item(211, 1).Value  = 3
$ws.Cells.Item(211, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(211, 3).Value  = "Coquimbo"
$ws.Cells.Item(211, 4).Value  = 44524
$ws.Cells.Item(211, 5).Value  = 5
$ws.Cells.Item(211, 6).Value  = 100112028
$ws.Cells.Item(211, 7).Value  = "Sandia"
$ws.Cells.Item(211, 8).Value  = "Sin especificar"
$ws.Cells.Item(211, 9).Value  = "Primera"
$ws.Cells.Item(211, 10).Value = 280
$ws.Cells.Item(211, 11).Value = 700
$ws.Cells.Item(211, 12).Value = 750
$ws.Cells.Item(211, 13).Value = 729
$ws.Cells.Item(211, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(211, 15).Value = "Perú"
$ws.Cells.Item(211, 16).Value = 729
$ws.Cells.Item(211, 17).Value = 1
$ws.Cells.Item(211, 18).Value = "Hortaliza"
